# fix JobDescription convertion datetime from string
# Add a new task row (SC-11) to the tracker and keep the table sorted by Priority.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right under the header, shifting the existing tasks down.
$ws.Rows("2:2").Insert()

# Insert() leaves stray empty cells behind in the columns that have no
# value on the new row (C and E) - drop them so the row matches the rest
# of the table (which only has cells for the columns actually in use).
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# New task entry.
$ws.Range("A2").Value = "SC-11"
$ws.Range("A2").Font.Size = 13
$ws.Range("A2").Font.Name = "Helvetica Neue"
$ws.Range("B2").Value = "Добавть юнит-тесты для JobDescription маппинга из JSON"
$ws.Range("D2").Value = "High"

# Re-apply (and persist) the table's "sort by Priority" rule now that the
# table grew by one row.
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("D2:D12"))
$sort.SetRange($ws.Range("A2:E12"))
$sort.Header = 0
$sort.Apply()

# The table is kept sorted by Priority (and, within a priority, by task id) -
# reorder the existing "Middle" priority rows to restore that order now that
# a new row was inserted above them.
$ws.Range("A9").Value = "SC-10"
$ws.Range("B9").Value = "Создать Task для записи файла"
$ws.Range("A10").Value = "SC-2"
$ws.Range("B10").Value = " Создать ТАСК удаление файла" + [char]160
$ws.Range("A11").Value = "SC-6"
$ws.Range("B11").Value = " Создать ТАСК для GET запросов (погоды)" + [char]160
$ws.Range("A12").Value = "SC-9"
$ws.Range("A12").Font.Size = 13
$ws.Range("A12").Font.Name = "Helvetica Neue"
$ws.Range("B12").Value = "Создать Task для чтения файла"
$ws.Range("D12").Value = "Middle"

# Column B needs to widen to fit the new, longer task description
# (mirrors the "best fit" width Excel recalculates for the longer text).
$ws.Columns("B:B").ColumnWidth = 50.25

# Restore the selection to where the user left it after editing.
$ws.Range("E8").Select()
